$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "29.916.30"
$ws.Range("E2").Value = "  +0.05%  "

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.876.11"
$ws.Range("E3").Value = "  -0.64%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5 (XRP) ---
$ws.Range("D5").Value = "'0.7434"
$ws.Range("E5").Value = "  -3.76%  "

# --- Row 6 (BNB) ---
$ws.Range("D6").Value = "'242.55"
$ws.Range("E6").Value = "  -0.12%  "

# --- Row 7 (USDC) ---
$ws.Range("E7").Value = "  +0.00%  "

# --- Row 8 (Cardano) ---
$ws.Range("D8").Value = "'0.3151"
$ws.Range("E8").Value = "  +1.29%  "

# --- Row 9 (Dogecoin) ---
$ws.Range("D9").Value = "'0.07240"
$ws.Range("E9").Value = "  +0.74%  "

# --- Row 10 (Solana) ---
$ws.Range("D10").Value = "'24.71"
$ws.Range("E10").Value = "  -3.53%  "

# --- Row 11 (TRON) ---
$ws.Range("D11").Value = "'0.08405"
$ws.Range("E11").Value = "  -2.25%  "

# --- Row 12 (Polygon(wasWrappedEther)) ---
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7522"
$ws.Range("E12").Value = "  -1.60%  "

# --- Row 13 (Polkadot(wasPolygon)) ---
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.431"
$ws.Range("E13").Value = "  +1.01%  "

# --- Row 14 (WrappedEther(wasPolkadot)) ---
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.875.66"
$ws.Range("E14").Value = "  -4.06%  "

# --- Row 15 (Litecoin) ---
$ws.Range("D15").Value = "'92.54"
$ws.Range("E15").Value = "  -1.29%  "

# --- Row 16 (WrappedBTC) ---
$ws.Range("D16").Value = "29.923.61"
$ws.Range("E16").Value = "  -0.14%  "

# --- Row 17 (Uniswap) ---
$ws.Range("D17").Value = "'6.094"
$ws.Range("E17").Value = "  -1.56%  "

# --- Row 18 (BitcoinCash) ---
$ws.Range("D18").Value = "'247.80"
$ws.Range("E18").Value = "  +1.45%  "

# --- Row 19 (Avalanche) ---
$ws.Range("E19").Value = "  -1.36%  "

# --- Row 20 (ShibaInu) ---
$ws.Range("D20").Value = "'0.000007859"
$ws.Range("E20").Value = "  +0.50%  "

# --- Row 21 (Dai) ---
$ws.Range("E21").Value = "  +0.14%  "

# --- Row 22 (WrappedliquidstakedEther2.0) ---
$ws.Range("D22").Value = "2.128.17"
$ws.Range("E22").Value = "  -3.52%  "

# --- Row 23 (Chainlink) ---
$ws.Range("D23").Value = "'8.054"
$ws.Range("E23").Value = "  +0.77%  "

# --- Row 24 (BinanceUSD) ---
$ws.Range("E24").Value = "  +0.02%  "

# --- Row 25 (Stellar) ---
$ws.Range("D25").Value = "'0.1562"
$ws.Range("E25").Value = "  -5.18%  "

# --- Row 27 (Monero) ---
$ws.Range("D27").Value = "'165.25"
$ws.Range("E27").Value = "  +2.10%  "

# --- Row 29 (LidoDAOToken) ---
$ws.Range("E29").Value = "  +0.25%  "

# --- Row 30 (Toncoin) ---
$ws.Range("D30").Value = "'1.516"
$ws.Range("E30").Value = "  +5.12%  "

# --- Row 31 (Filecoin) ---
$ws.Range("E31").Value = "  +1.77%  "

# --- Row 32 (PancakeSwap) ---
$ws.Range("D32").Value = "'1.537"
$ws.Range("E32").Value = "  +0.30%  "

# --- Row 33 (InternetComputer(DFINITY)) ---
$ws.Range("D33").Value = "'4.287"
$ws.Range("E33").Value = "  +4.48%  "

# --- Row 34 (Hedera) ---
$ws.Range("D34").Value = "'0.05347"
$ws.Range("E34").Value = "  -1.69%  "

# --- Row 35 (ARBITRUM) ---
$ws.Range("D35").Value = "'1.238"
$ws.Range("E35").Value = "  -0.19%  "

# --- Row 36 (ImmutableX) ---
$ws.Range("D36").Value = "'0.7526"
$ws.Range("E36").Value = "  +0.88%  "

# --- Row 37 (Frax) ---
$ws.Range("D37").Value = "'0.9996"
$ws.Range("E37").Value = "  -0.37%  "

# --- Row 38 (HuobiToken) ---
$ws.Range("D38").Value = "'2.690"
$ws.Range("E38").Value = "  -0.16%  "

# --- Row 39 (VeChain) ---
$ws.Range("E39").Value = "  -0.02%  "

# --- Row 40 (MXToken) ---
$ws.Range("D40").Value = "'2.754"
$ws.Range("E40").Value = "  -0.94%  "

# --- Row 41 (TheSandbox) ---
$ws.Range("D41").Value = "'0.4547"
$ws.Range("E41").Value = "  +1.83%  "

# --- Row 42 (Maker) ---
$ws.Range("D42").Value = "1.112.64"
$ws.Range("E42").Value = "  +0.22%  "

# --- Row 43 (FraxShare) ---
$ws.Range("D43").Value = "'6.058"
$ws.Range("E43").Value = "  -0.50%  "

# --- Row 44 (Aave) ---
$ws.Range("D44").Value = "'72.56"
$ws.Range("E44").Value = "  -1.12%  "

# --- Row 45 (TrustWalletToken) ---
$ws.Range("D45").Value = "'0.8568"
$ws.Range("E45").Value = "  +0.78%  "

# --- Row 46 (PaxDollar) ---
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.08%  "

# --- Row 47 (Quant) ---
$ws.Range("D47").Value = "'103.39"
$ws.Range("E47").Value = "  +0.44%  "

# --- Row 48 (RenderToken(wasAptos)) ---
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.859"
$ws.Range("E48").Value = "  -0.67%  "

# --- Row 49 (Aptos(wasRenderToken)) ---
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.625"
$ws.Range("E49").Value = "  -0.03%  "

# --- Row 50 (RocketPoolETH) ---
$ws.Range("D50").Value = "2.025.81"
$ws.Range("E50").Value = "  -2.44%  "

# --- Row 51 (SynthetixNetwork) ---
$ws.Range("D51").Value = "'2.909"
$ws.Range("E51").Value = "  -2.53%  "

